$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.886.70"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.993.46"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.87%  "
$ws.Range("E9").Value = "  -4.52%  "
$ws.Range("E10").Value = "  -5.30%  "
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("D12").Value = "2.288.70"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.789"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.24%  "
$ws.Range("D17").Value = "1.996.64"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "36.818.28"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "235.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.123"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.72%  "
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("E33").Value = "  -5.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("E35").Value = "  -9.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.447.46"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.99%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0208"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.39%  "
$ws.Range("D50").Value = "2.181.03"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.14%  "
